$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-339) from 2023-09-10 to 2023-09-11, keeping the existing
# date number format/style intact.
$ws.Range("C2:C339").Value = "2023-09-11"
